$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'257.35"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'0.13%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'26.97"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-1.82%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'4.682"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-10.06%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.05885"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-0.47%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'6.647"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-0.30%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.8576"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-0.74%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.9546"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-4.81%"
$ws.Range("E8").Style = "Normal"
$ws.Range("B9").Value = 'WazirX'
$ws.Range("C9").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D9").Value = "'0.1408"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-0.71%"
$ws.Range("E9").Style = "Normal"
$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D10").Value = "'0.03934"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'9.81%"
$ws.Range("E10").Style = "Normal"
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").Value = "'0.07094"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-1.27%"
$ws.Range("E11").Style = "Normal"
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").Value = "'0.03183"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'1.10%"
$ws.Range("E12").Style = "Normal"
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").Value = "'0.09169"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.64%"
$ws.Range("E13").Style = "Normal"
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").Value = "'0.001549"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'0.48%"
$ws.Range("E14").Style = "Normal"
$ws.Range("B15").Value = 'One'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D15").Value = "'0.0006031"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-0.81%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.006204"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'5.06%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.516"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'0.74%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.203"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-1.97%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'2.227"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'0.15%"
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'-2.24%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.1293"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-1.05%"
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'8.80%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04238"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'1.26%"
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'0.40%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004298"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-4.76%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0001200"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'0.02%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.0001937"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'-0.05%"
$ws.Range("E27").Style = "Normal"
$ws.Range("E40").Value = "'0.35%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.006291"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'12.15%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1104"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'0.28%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002440"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'35.59%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.01143"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'7.04%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005441"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'0.26%"
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'0.02%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.06000"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-45.00%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.1675"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'7,592.92%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.00002100"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'0.02%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.0002000"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.02%"
$ws.Range("E50").Style = "Normal"
